# LOM3083.docx edit: rotate the discipline-description content blocks.
#
# The edit reshuffles which "body" paragraph holds which block of text
# (the paragraph *shells* - styles/labels - stay exactly where they are):
#
#   Objetivos body            <- old Programa-resumido body
#   Docente(s) body           <- old Objetivos body
#   Programa-resumido body    <- old Programa (numbered list) body
#   Programa body             <- old "Método:" content (Avaliação block)
#   Avaliação/"Método:"       <- old "Critério:" content
#   Avaliação/"Critério:"     <- old "Norma de recuperação:" content
#   Avaliação/"Norma de..."   <- old Bibliografia body
#   Bibliografia body         <- old Docente(s) body
#
# We capture every original value first (so later writes never race
# against not-yet-read text), then write the rotated values back.

$d = $word.ActiveDocument

function Get-ParaText($index) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range.Duplicate
    $r.End = $r.End - 1   # drop the trailing paragraph mark
    return $r.Text
}

function Set-ParaText($index, $value) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range.Duplicate
    $r.End = $r.End - 1   # keep the paragraph mark / pPr untouched
    $r.Text = $value
}

# --- 1. Capture all current values -----------------------------------
$objetivosBody     = Get-ParaText 6    # "Como parte fundamental..."
$docenteBody       = Get-ParaText 8    # "1176388 - Luiz Tadeu Fernandes Eleno"
$progResumidoBody  = Get-ParaText 10   # "Introdução à transferência de calor..."
$programaBody      = Get-ParaText 12   # "1.Transferência de calor por condução... uphill."
$bibliografiaBody  = Get-ParaText 16   # "INCROPERA ... 2003." (with embedded vertical-tab breaks)

# The "Avaliação" paragraph (index 14) mixes three bold labels with
# three plain-text contents; grab each content piece individually via
# Find so the bold "Método:"/"Critério:"/"Norma de recuperação:" runs
# are left completely alone.
$avalPara = $d.Paragraphs.Item(14).Range

$metodoRange = $avalPara.Duplicate
$metodoRange.Find.Execute("Aulas expositivas e interação em grupo para a solução de problemas.") | Out-Null
$metodoContent = $metodoRange.Text

$criterioRange = $avalPara.Duplicate
$criterioRange.Find.Execute("Média aritmética (M) provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,4*P2+0.3*T") | Out-Null
$criterioContent = $criterioRange.Text

$normaRange = $avalPara.Duplicate
$normaRange.Find.Execute("Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação") | Out-Null
$normaContent = $normaRange.Text

# --- 2. Write the rotated values back ---------------------------------

# Objetivos body <- old Programa-resumido body
Set-ParaText 6 $progResumidoBody

# Docente(s) body <- old Objetivos body
Set-ParaText 8 $objetivosBody

# Programa-resumido body <- old Programa (numbered list) body
Set-ParaText 10 $programaBody

# Programa body <- old "Método:" content
Set-ParaText 12 $metodoContent

# Avaliação block: shift Método <- Critério <- Norma <- Bibliografia.
# Write back-to-front (Norma, then Critério, then Método) so each
# Find.Execute below still locates the *original*, not-yet-overwritten
# text - otherwise an earlier write could shadow a later search target.
$normaRange2 = $d.Paragraphs.Item(14).Range.Duplicate
$normaRange2.Find.Execute("Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação") | Out-Null
$normaRange2.Text = $bibliografiaBody

$criterioRange2 = $d.Paragraphs.Item(14).Range.Duplicate
$criterioRange2.Find.Execute("Média aritmética (M) provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,4*P2+0.3*T") | Out-Null
$criterioRange2.Text = $normaContent

$metodoRange2 = $d.Paragraphs.Item(14).Range.Duplicate
$metodoRange2.Find.Execute("Aulas expositivas e interação em grupo para a solução de problemas.") | Out-Null
$metodoRange2.Text = $criterioContent

# Bibliografia body <- old Docente(s) body
Set-ParaText 16 $docenteBody

Write-Output "done"
